$d = $word.ActiveDocument

# Locate the target sentence (the run to be split).
$found = $d.Content
$found.Find.Execute("; so, it gives a best as possible fingerprint protection.")
$start = $found.Start

# Isolate the first 3 characters ("; s") of that run into their own run by
# toggling a formatting property first -- this forces the engine to split
# the run without merging/restamping rsid on the untouched remainder.
$r1 = $d.Range($start, $start + 3)
$r1.Bold = 1

# Now replace just that isolated run's text: "; s" -> ". S"
$r1.Text = ". S"

# Restore formatting so the new run's rPr matches its Cambria-only sibling.
$r1.Bold = 0
